{"js": "const body = context.document.body;\nconst doc = context.document;\n\n// 1) \"This exploratory study ...\" -> \"This innovative study ...\"\nconst exploratoryResults = body.search(\"exploratory\", { matchCase: true });\nexploratoryResults.load(\"text\");\nawait context.sync();\nif (exploratoryResults.items.length !== 1) {\n  throw new Error(\"expected exactly one match for 'exploratory', got \" + exploratoryResults.items.length);\n}\nexploratoryResults.items[0].insertText(\"innovative\", \"Replace\");\nawait context.sync();\n\n// 2) \"Reviewers?\" -> \"Reviewers\", and move the _GoBack bookmark from the end\n//    of the document to right after this heading.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nlet reviewersParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Reviewers?\") !== -1) {\n    reviewersParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!reviewersParagraph) {\n  throw new Error(\"could not find the 'Reviewers?' paragraph\");\n}\nconst reviewersParaRange = reviewersParagraph.getRange();\nconst reviewersOoxml = reviewersParaRange.getOoxml();\nawait context.sync();\nlet reviewersXml = reviewersOoxml.value;\nif (reviewersXml.indexOf(\"Reviewers?\") === -1) {\n  throw new Error(\"unexpected paragraph contents while editing 'Reviewers?'\");\n}\nreviewersXml = reviewersXml.replace(\"Reviewers?\", \"Reviewers\");\nreviewersParaRange.insertOoxml(reviewersXml, \"Replace\");\nawait context.sync();\n\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst reviewersHeading = body.search(\"Reviewers\", { matchCase: true });\nreviewersHeading.load(\"text\");\nawait context.sync();\nif (reviewersHeading.items.length !== 1) {\n  throw new Error(\"expected exactly one match for 'Reviewers', got \" + reviewersHeading.items.length);\n}\nconst reviewersEnd = reviewersHeading.items[0].getRange(\"End\");\nreviewersEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Author contributions: add \"D. Corina, \" to the list of authors who\n//    developed the study concept.\nconst conceptResults = body.search(\"V.A. Marchman, and A. Fernald developed\", { matchCase: true });\nconceptResults.load(\"text\");\nawait context.sync();\nif (conceptResults.items.length !== 1) {\n  throw new Error(\"expected exactly one match for concept sentence, got \" + conceptResults.items.length);\n}\nconceptResults.items[0].insertText(\"V.A. Marchman, D. Corina, and A. Fernald developed\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"This exploratory study ...\" -> \"This innovative study ...\"\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Text = \"exploratory\"\n$found1 = $find1.Execute()\nif (-not $found1) {\n    throw \"Could not find 'exploratory'\"\n}\n$range1.Text = \"innovative\"\n\n# 2) \"Reviewers?\" -> \"Reviewers\", and move the _GoBack bookmark from the end\n#    of the document to right after this heading.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"Reviewers?\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n    throw \"Could not find 'Reviewers?'\"\n}\n# $range2 now spans \"Reviewers?\" -- bookmark just the trailing \"?\" so that,\n# once we clear its text, the bookmark collapses to a single point right\n# after \"Reviewers\" (Bookmarks.Add with an already-existing name moves it).\n$qMark = $d.Range($range2.End - 1, $range2.End)\n$d.Bookmarks.Add(\"_GoBack\", $qMark)\n$qMark.Text = \"\"\n\n# 3) Author contributions: add \"D. Corina, \" to the list of authors who\n#    developed the study concept.\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.Text = \"V.A. Marchman, and A. Fernald developed\"\n$found3 = $find3.Execute()\nif (-not $found3) {\n    throw \"Could not find the author-contributions sentence\"\n}\n$range3.Text = \"V.A. Marchman, D. Corina, and A. Fernald developed\"\n"}
